$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 169, shifting existing rows 169:242 down to 170:243
$ws.Rows("169:169").Insert()

# Populate the newly inserted row 169 with the new record's data
$ws.Range("A169").Value = 3
$ws.Range("B169").Value = "Femacal de La Calera"
$ws.Range("C169").Value = "Coquimbo"
$ws.Range("D169").Value = 44510
$ws.Range("E169").Value = 5
$ws.Range("F169").Value = 100112031
$ws.Range("G169").Value = "Poroto verde"
$ws.Range("H169").Value = "Magnum"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 70
$ws.Range("K169").Value = 37000
$ws.Range("L169").Value = 38000
$ws.Range("M169").Value = 37500
$ws.Range("N169").Value = "`$/malla 25 kilos"
$ws.Range("O169").Value = "Provincia de Limarí"
$ws.Range("P169").Value = 1500
$ws.Range("Q169").Value = 25
$ws.Range("R169").Value = "Hortaliza"
